$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update row 2 (Monica Francisca Gastanbide -> Max Verstapen Quimera)
$ws.Range("A2").Value = "Max"
$ws.Range("B2").Value = "Verstapen"
$ws.Range("C2").Value = "Quimera"
$ws.Range("D2").Value = 202102001
$ws.Range("E2").Value = "MaxF1V"

# Update row 3 (Charles Antonio Lecrec Montez -> Luis Alberto Hamilton Vera)
$ws.Range("A3").Value = "Luis"
$ws.Range("B3").Value = "Alberto"
$ws.Range("C3").Value = "Hamilton Vera"
$ws.Range("D3").Value = 202102002
$ws.Range("E3").Value = "LHalF1"

# Widen column C to fit new content (target stored width 16.42578125 chars)
$ws.Columns.Item(3).ColumnWidth = 15.65
